$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1986.3158
$ws.Range("I98").Value = 1583.75
$ws.Range("J98").Value = 4133.3335
$ws.Range("K98").Value = 1583.75
$ws.Range("L98").Value = 4133.3335
$ws.Range("M98").Value = -85.75
$ws.Range("N98").Value = -7129.3335
$ws.Range("H112").Value = 2656.3635
$ws.Range("J112").Value = 3340.625
$ws.Range("L112").Value = 10021.875
$ws.Range("N112").Value = -12237.875
$ws.Range("H122").Value = 1986.3158
$ws.Range("I122").Value = 1583.75
$ws.Range("J122").Value = 4133.3335
$ws.Range("K122").Value = 4751.25
$ws.Range("L122").Value = 12400.0005
$ws.Range("M122").Value = -2301.25
$ws.Range("N122").Value = -17300.0005
$ws.Range("H137").Value = 7464225
$ws.Range("I137").Value = 11629491
$ws.Range("J137").Value = 1455.8334
$ws.Range("K137").Value = 34888473
$ws.Range("L137").Value = 4367.5002
$ws.Range("M137").Value = -34885923
$ws.Range("N137").Value = -9467.5002
$ws.Range("H138").Value = 1999.8
$ws.Range("I138").Value = 2032.1
$ws.Range("J138").Value = 1981.3429
$ws.Range("K138").Value = 6096.299999999999
$ws.Range("L138").Value = 5944.028700000001
$ws.Range("M138").Value = -956.2999999999993
$ws.Range("N138").Value = -16224.0287

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3500.17
$ws.Range("I32").Value = 3531.3586
$ws.Range("J32").Value = 3141.5
$ws.Range("K32").Value = 3531.3586
$ws.Range("L32").Value = 3141.5
$ws.Range("M32").Value = -3244.3586
$ws.Range("N32").Value = -3715.5
$ws.Range("H122").Value = 1394.305
$ws.Range("I122").Value = 990.84
$ws.Range("K122").Value = 2972.52
$ws.Range("M122").Value = -522.52
$ws.Range("H132").Value = 2196193
$ws.Range("I132").Value = 1293.3422
$ws.Range("J132").Value = 6585992.5
$ws.Range("K132").Value = 3880.0266
$ws.Range("L132").Value = 19757977.5
$ws.Range("M132").Value = -1350.0266
$ws.Range("N132").Value = -19763037.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3565.4531
$ws.Range("I134").Value = 1467.7435
$ws.Range("J134").Value = 6837.88
$ws.Range("K134").Value = 4403.2305
$ws.Range("L134").Value = 20513.64
$ws.Range("M134").Value = -1868.2305
$ws.Range("N134").Value = -25583.64

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 30067.75
$ws.Range("J64").Value = 30067.75
$ws.Range("L64").Value = 30067.75
$ws.Range("N64").Value = -30563.75
$ws.Range("H67").Value = 30067.75
$ws.Range("J67").Value = 30067.75
$ws.Range("L67").Value = 30067.75
$ws.Range("N67").Value = -31783.75
$ws.Range("H94").Value = 1335
$ws.Range("I94").Value = 1212
$ws.Range("K94").Value = 1212
$ws.Range("M94").Value = -761
$ws.Range("H132").Value = 3008.5386
$ws.Range("I132").Value = 2618.5715
$ws.Range("J132").Value = 3152.2104
$ws.Range("K132").Value = 7855.7145
$ws.Range("L132").Value = 9456.6312
$ws.Range("M132").Value = -5325.7145
$ws.Range("N132").Value = -14516.6312
$ws.Range("H134").Value = 3189.5454
$ws.Range("I134").Value = 1286.2
$ws.Range("J134").Value = 4775.6665
$ws.Range("K134").Value = 3858.6
$ws.Range("L134").Value = 14326.9995
$ws.Range("M134").Value = -1323.6
$ws.Range("N134").Value = -19396.9995

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2600
$ws.Range("I80").Value = 1400
$ws.Range("J80").Value = 2709.0908
$ws.Range("K80").Value = 4200
$ws.Range("L80").Value = 8127.2724
$ws.Range("M80").Value = -3264
$ws.Range("N80").Value = -9999.2724
$ws.Range("H83").Value = 2600
$ws.Range("I83").Value = 1400
$ws.Range("J83").Value = 2709.0908
$ws.Range("K83").Value = 12600
$ws.Range("L83").Value = 24381.8172
$ws.Range("M83").Value = -7920
$ws.Range("N83").Value = -33741.8172

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6759.9
$ws.Range("I70").Value = 5299.857
$ws.Range("K70").Value = 5299.857
$ws.Range("M70").Value = -5029.857
$ws.Range("H73").Value = 6759.9
$ws.Range("I73").Value = 5299.857
$ws.Range("K73").Value = 5299.857
$ws.Range("M73").Value = -4363.857

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("H94").Value = 24110
$ws.Range("J94").Value = 24110
$ws.Range("L94").Value = 24110
$ws.Range("N94").Value = -25462
$ws.Range("H122").Value = 9516.235000000001
$ws.Range("I122").Value = 11190.546
$ws.Range("J122").Value = 6446.6665
$ws.Range("K122").Value = 33571.638
$ws.Range("L122").Value = 19339.9995
$ws.Range("M122").Value = -31121.638
$ws.Range("N122").Value = -24239.9995
$ws.Range("H132").Value = 26318448
$ws.Range("I132").Value = 30305596
$ws.Range("J132").Value = 3279
$ws.Range("K132").Value = 90916788
$ws.Range("L132").Value = 9837
$ws.Range("M132").Value = -90914258
$ws.Range("N132").Value = -14897
$ws.Range("H136").Value = 8475338
$ws.Range("I136").Value = 10204733
$ws.Range("J136").Value = 1303.4
$ws.Range("K136").Value = 30614199
$ws.Range("L136").Value = 3910.2
$ws.Range("M136").Value = -30611649
$ws.Range("N136").Value = -9010.200000000001
$ws.Range("M81").ClearContents()
$ws.Range("M84").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31248
$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96240
$ws.Range("H113").Value = 898.069
$ws.Range("J113").Value = 1034.0454
$ws.Range("L113").Value = 3102.1362
$ws.Range("N113").Value = -7442.1362
$ws.Range("H132").Value = 2349.3513
$ws.Range("I132").Value = 2023.6111
$ws.Range("J132").Value = 2657.9473
$ws.Range("K132").Value = 6070.8333
$ws.Range("L132").Value = 7973.841899999999
$ws.Range("M132").Value = -3540.8333
$ws.Range("N132").Value = -13033.8419
$ws.Range("H136").Value = 7821715.5
$ws.Range("I136").Value = 10428116
$ws.Range("J136").Value = 2513.125
$ws.Range("K136").Value = 31284348
$ws.Range("L136").Value = 7539.375
$ws.Range("M136").Value = -31281798
$ws.Range("N136").Value = -12639.375
